# fix SH688098 and SH600636 conflict
# Append 4 new date/value rows (120-123) to the end of the price-history sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates  = @("20201123", "20201207", "20201228", "20210105")
$values = @(915.0, 1087.0, 799.0, 927.0)

$startRow = 120
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    # Leading apostrophe forces the numeric-looking date string to be
    # stored as literal text, matching the existing column-A cells.
    $ws.Cells.Item($row, 1).Value = "'" + $dates[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
